# Add new columns I ("I0") and J ("IF") with per-row data, rows 1-81.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - match style of other header cells (bold/centered/bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$I = @(7,7,12,9,9,7,7,8,7,7,8,8,9,8,8,7,9,8,7,7,8,9,8,8,7,8,8,9,8,8,8,8,8,8,9,8,9,10,8,9,7,8,8,7,8,8,8,8,8,8,7,8,7,9,7,7,7,8,8,7,9,9,9,8,9,8,8,8,8,8,8,9,5,8,5,5,5,4,6,2)
$J = @(7,8,12,9,9,8,8,8,8,7,8,8,9,8,8,7,9,8,7,7,8,9,8,8,7,8,8,9,8,8,8,8,8,8,9,8,9,10,8,9,7,8,8,7,8,8,8,8,8,8,7,8,7,9,8,7,7,8,8,7,9,9,9,8,9,8,8,8,8,8,8,9,6,8,5,5,5,4,6,2)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}
